$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh (prices + 1h volume %); row 10/11 also swap
# Toncoin <-> Dogecoin ranking order.
#
# The Price column (D) holds values that look numeric ("63.444.17",
# "1.00", ...) but are stored as TEXT in the workbook. Setting a cell's
# .Value with such a string makes Excel auto-coerce it to a real number,
# so we force a text number format on each Price cell right before writing
# it, to keep it text exactly like the original inlineStr cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.444.17'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.078.04'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.52'
$ws.Range("E5").Value = '  -1.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.75'
$ws.Range("E6").Value = '  +1.20%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.072.82'
$ws.Range("E8").Value = '  -0.56%  '

$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  -1.66%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("E12").Value = '  -3.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000223'
$ws.Range("E13").Value = '  +2.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.62'
$ws.Range("E14").Value = '  -2.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.581.53'
$ws.Range("E15").Value = '  -0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.536.61'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("E17").Value = '  +0.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.078.45'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("E19").Value = '  -2.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '472.80'
$ws.Range("E20").Value = '  -3.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.38'
$ws.Range("E21").Value = '  -2.24%  '

$ws.Range("E22").Value = '  -2.63%  '

$ws.Range("E23").Value = '  -2.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.68'
$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.20'
$ws.Range("E25").Value = '  -1.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.99'
$ws.Range("E28").Value = '  -5.91%  '

$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.16'
$ws.Range("E30").Value = '  -1.75%  '

$ws.Range("E31").Value = '  -4.59%  '

$ws.Range("E32").Value = '  +1.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '57.33'
$ws.Range("E33").Value = '  -2.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.32'
$ws.Range("E34").Value = '  -7.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  +4.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '489.02'
$ws.Range("E36").Value = '  -6.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.98'
$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.231.03'
$ws.Range("E38").Value = '  +2.85%  '

$ws.Range("E39").Value = '  -1.87%  '

$ws.Range("E40").Value = '  -1.21%  '

$ws.Range("E41").Value = '  -0.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.09'
$ws.Range("E42").Value = '  -1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.63'
$ws.Range("E43").Value = '  -1.97%  '

$ws.Range("E44").Value = '  -2.11%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '123.50'
$ws.Range("E46").Value = '  +1.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.05'
$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.03'
$ws.Range("E48").Value = '  -3.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0529'
$ws.Range("E49").Value = '  +4.68%  '

$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  +4.05%  '
